# [FINALIZAR RUTINA] End routine button functionality added
#
# Updates the tracker model sheet:
#  - Renames the "push_increment_units / pull_increment / isometric_increment"
#    header trio (I3:K3) to the new "*_time_increment" names, and switches
#    their sample-row units from "Interger" to "Time (seconds)".
#  - Drops the now-unused "Time (seconds)" / "Time (minutes)" sample cells
#    in I5:K6.
#  - Replaces the old routine-group placeholders ("circuit", "anaerobic",
#    "time") in D12:D14 with the generic "[exercise_plan_id]" placeholder.
#  - Adds a new "exercise_plan_end" sample row (row 15) mirroring the
#    existing "exercise_plan_start" rows.
#  - Adds a small "Input JSON" legend table (rows 20-21) documenting the
#    new exercises_summary field, with highlighted header/key cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 3 headers: the push/pull/isometric "units" columns now describe a
# time increment instead of a generic increment.
# ---------------------------------------------------------------------
$ws.Range("I3").Value = "push_time_increment"
$ws.Range("J3").Value = "pull_time_increment"
$ws.Range("K3").Value = "isometric_time_increment"

# Row 4 sample values switch from "Interger" to "Time (seconds)".
$ws.Range("I4:K4").Value = "Time (seconds)"

# Rows 5 & 6 no longer carry sample values in I:K.
$ws.Range("I5:K5").ClearContents()
$ws.Range("I6:K6").ClearContents()

# ---------------------------------------------------------------------
# Rows 12-14: routine-group placeholders become the generic
# "[exercise_plan_id]" value.
# ---------------------------------------------------------------------
$ws.Range("D12").Value = "[exercise_plan_id]"
$ws.Range("D13").Value = "[exercise_plan_id]"
$ws.Range("D14").Value = "[exercise_plan_id]"

# ---------------------------------------------------------------------
# New row 15: exercise_plan_end, mirroring the exercise_plan_start rows
# above it (same formatting, pulled from row 14).
# ---------------------------------------------------------------------
$ws.Range("B14:K14").Copy()
$ws.Range("B15:K15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B15").Value = "YYYY-MM-DD"
$ws.Range("C15").Value = "exercise_plan_end"
$ws.Range("D15").Value = "[exercise_plan_id]"
$ws.Range("E15:K15").Value = "None"

# ---------------------------------------------------------------------
# New "Input JSON" mini legend (rows 20-21) describing the
# exercises_summary field used by the end-routine button.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = "field name"
$ws.Range("D20").Value = "field description"
$ws.Range("B21").Value = "Input JSON"
$ws.Range("C21").Value = "exercises_summary"
$ws.Range("D21").Value = '{routine_group: chest, exercise_reps: {"push-ups": 20, "muscle-up":5}}'

# Highlight fills for the new legend row.
$ws.Range("C20:D20").Interior.Color = 12566463
$ws.Range("B21").Interior.Color = 9555625
$ws.Range("C21").Interior.Color = 10086143

# ---------------------------------------------------------------------
# Column width tweaks.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws.Columns.Item(10).ColumnWidth = 20
$ws.Columns.Item(11).ColumnWidth = 23.5
$ws.Columns.Item(12).ColumnWidth = 24.5

# Selection, matching where the author last left the cursor.
$ws.Range("E14").Select()
